$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1) "vaccines links" sheet: new row 8 (Wikipedia / SARS-CoV-2-Impfstoff) ---
# Entered first so shared-string allocation order matches the source edit.
$ws2.Range("A8").Value = "Wikipedia: "
$ws2.Range("B8").Value = "COVID-19 vaccine"
$ws2.Range("C8").Value = "SARS-CoV-2-Impfstoff"
$ws2.Hyperlinks.Add($ws2.Range("D8"), "https://de.wikipedia.org/wiki/SARS-CoV-2-Impfstoff")

# --- 2) "vaccines" sheet: reserved doses for Moderna (N3) updated ---
$ws1.Range("N3").Value = 13500000

# --- 3) "vaccines" sheet: new row 6 - Curevac / Zorecimeran ---
$ws1.Range("A6").Value = "Curevac"
$ws1.Range("B6").Value = "Zorecimeran"
$ws1.Hyperlinks.Add($ws1.Range("C6"), "https://de.wikipedia.org/wiki/CVnCoV")
$ws1.Range("D6").Value = "DE"
$ws1.Range("E6").Value = "mRNA vaccine"
$ws1.Range("F6").Value = "https://www.swissmedic.ch/swissmedic/de/home/ueber-uns/publikationen/video/different-types-of-vaccine.html"
$ws1.Range("G6").NumberFormat = $ws1.Range("G4").NumberFormat
$ws1.Range("G6").Value = "TBD"
$ws1.Range("H6").Value = 2
$ws1.Range("I6").Value = "-"
$ws1.Range("J6").Value = "not approved"
$ws1.Range("K6").Value = "not known"
$ws1.Range("L6").Value = "undergoing approval process"
$ws1.Range("M6").Value = "TBD"
$ws1.Range("N6").Value = 5000000
$ws1.Range("O6").NumberFormat = $ws1.Range("O4").NumberFormat
$ws1.Range("O6").Value = 44230

# --- 4) "vaccines" sheet: new row 7 - Novavax / NVX-CoV2373 ---
$ws1.Range("A7").Value = "Novavax"
$ws1.Range("B7").Value = "NVX-CoV2373"
$ws1.Range("C7").Value = "https://en.wikipedia.org/wiki/Novavax_COVID-19_vaccine"
$ws1.Range("D7").Value = "US"
$ws1.Range("E7").Value = "recombinant nanoparticle vaccine"
$ws1.Range("G7").NumberFormat = $ws1.Range("G4").NumberFormat
$ws1.Range("G7").Value = "TBD"
$ws1.Range("H7").Value = 2
$ws1.Range("I7").Value = "-"
$ws1.Range("J7").Value = "not approved"
$ws1.Range("K7").Value = "regular fridge temperature"
$ws1.Range("L7").Value = "undergoing approval process"
$ws1.Range("M7").Value = "TBD"
$ws1.Range("N7").Value = 6000000
$ws1.Range("O7").NumberFormat = $ws1.Range("O4").NumberFormat
$ws1.Range("O7").Value = 44230

# --- 5) view state: "vaccines" tab now active/selected, "vaccines links" no longer ---
$ws2.Activate()
$ws2.Range("A8:D8").Select()
$excel.ActiveWindow.Zoom = 85

$ws1.Activate()
$ws1.Range("K14").Select()
